$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.974.85"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "2.353.35"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "'239.08"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'73.70"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.594"
$ws.Range("E9").Value = "  +8.35%  "
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("D11").Value = "'57.32"
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "'32.41"
$ws.Range("E12").Value = "  +9.89%  "
$ws.Range("D13").Value = "'7.28"
$ws.Range("E13").Value = "  +8.24%  "
$ws.Range("D14").Value = "'0.108"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").Value = "2.703.88"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "'16.55"
$ws.Range("E16").Value = "  -1.95%  "
$ws.Range("D17").Value = "'0.897"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").Value = "2.353.39"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("D19").Value = "43.864.79"
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").Value = "'6.74"
$ws.Range("E21").Value = "  +4.32%  "
$ws.Range("D22").Value = "'76.84"
$ws.Range("E22").Value = "  -1.68%  "
$ws.Range("D23").Value = "'258.91"
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("D24").Value = "'1.94"
$ws.Range("E24").Value = "  +22.29%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "'3.65"
$ws.Range("E26").Value = "  -2.81%  "
$ws.Range("E27").Value = "  -1.99%  "
$ws.Range("D28").Value = "'10.72"
$ws.Range("E28").Value = "  +2.00%  "
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").Value = "'22.59"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("D31").Value = "'175.35"
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("E32").Value = "  -3.47%  "
$ws.Range("E33").Value = "  +2.22%  "
$ws.Range("D34").Value = "'0.0761"
$ws.Range("E34").Value = "  +2.38%  "
$ws.Range("D35").Value = "'5.21"
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("D36").Value = "'5.47"
$ws.Range("E36").Value = "  +4.33%  "
$ws.Range("D37").Value = "'3.74"
$ws.Range("E37").Value = "  -4.93%  "
$ws.Range("E38").Value = "  -3.93%  "
$ws.Range("E39").Value = "  -3.30%  "
$ws.Range("D40").Value = "'0.0276"
$ws.Range("E40").Value = "  +1.58%  "
$ws.Range("E41").Value = "  +12.56%  "
$ws.Range("D42").Value = "'0.206"
$ws.Range("E42").Value = "  +12.78%  "
$ws.Range("D43").Value = "'18.89"
$ws.Range("E43").Value = "  -4.01%  "
$ws.Range("D44").Value = "'8.97"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").Value = "'4.70"
$ws.Range("E46").Value = "  +4.93%  "
$ws.Range("E47").Value = "  +6.22%  "
$ws.Range("D48").Value = "'57.77"
$ws.Range("E48").Value = "  +9.14%  "
$ws.Range("E49").Value = "  -1.09%  "
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("D51").Value = "'99.82"
$ws.Range("E51").Value = "  +1.15%  "
